$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 143 (weekly update: new price
# observations for Americana (o) Primera / Segunda on 2021-09-09 / serial 44448),
# pushing the previously-existing rows 143:149 down to 145:151.
$ws.Rows("143:144").Insert()

# New row 143: Ají, Americana (o), Primera
$ws.Cells.Item(143, 1).Value = 2
$ws.Cells.Item(143, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44448
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = 100112021
$ws.Cells.Item(143, 7).Value = "Ají"
$ws.Cells.Item(143, 8).Value = "Americana (o)"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 160
$ws.Cells.Item(143, 11).Value = 65000
$ws.Cells.Item(143, 12).Value = 70000
$ws.Cells.Item(143, 13).Value = 67500
$ws.Cells.Item(143, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(143, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(143, 16).Value = 2700
$ws.Cells.Item(143, 17).Value = 25
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# New row 144: Ají, Americana (o), Segunda
$ws.Cells.Item(144, 1).Value = 2
$ws.Cells.Item(144, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44448
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 100112021
$ws.Cells.Item(144, 7).Value = "Ají"
$ws.Cells.Item(144, 8).Value = "Americana (o)"
$ws.Cells.Item(144, 9).Value = "Segunda"
$ws.Cells.Item(144, 10).Value = 100
$ws.Cells.Item(144, 11).Value = 55000
$ws.Cells.Item(144, 12).Value = 60000
$ws.Cells.Item(144, 13).Value = 57500
$ws.Cells.Item(144, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(144, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(144, 16).Value = 2300
$ws.Cells.Item(144, 17).Value = 25
$ws.Cells.Item(144, 18).Value = "Hortaliza"
